$d = $word.ActiveDocument

function Replace-ParaContent($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $xmlFrag = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xmlFrag)
}

# --- Change 1: "(kapcsolattartó)" -> "(Kapcsolattartó)", split into two runs ---
# Locate the run " (kapcsolattartó)" that follows "- Lovas István"
$rng1 = $d.Range(0, $d.Content.End)
$rng1.Find.Execute(" (kapcsolattartó)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$runStart1 = $rng1.Start
$runEnd1 = $rng1.End

# Capitalise the "k" (still a single run at this point, safe to use plain Text assignment)
$kPos1 = $runStart1 + 2
$kRng1 = $d.Range($kPos1, $kPos1 + 1)
$kRng1.Text = "K"

# Split the run after "K": everything from just after "K" to the end of the run becomes a
# separate run, while " (K" stays behind untouched as its own run.
$tailRng1 = $d.Range($kPos1 + 1, $runEnd1)
$xmlTail1 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>apcsolattartó)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRng1.InsertXML($xmlTail1)

# --- Change 2: "(demonstrátor)" -> "(Demonstrátor)", split into two runs ---
$rng2 = $d.Range(0, $d.Content.End)
$rng2.Find.Execute(" (demonstrátor)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$runStart2 = $rng2.Start
$runEnd2 = $rng2.End

$dPos2 = $runStart2 + 2
$dRng2 = $d.Range($dPos2, $dPos2 + 1)
$dRng2.Text = "D"

$tailRng2 = $d.Range($dPos2 + 1, $runEnd2)
$xmlTail2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>emonstrátor)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRng2.InsertXML($xmlTail2)

# --- Changes 3-5: append " Lovas István, Balla Gábor" to the three bullet paragraphs,
#     and move the _GoBack bookmark out of the "csoportvezető" paragraph into the
#     "ifjúsági vezető" paragraph (at its very end). ---

$bulletPPr = '<w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="1418"/></w:tabs></w:pPr>'

# Find the three bullet paragraphs by their (stable) text prefixes.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("A t") -and $txt.Contains("funkciók osztálydiagramjai")) {
        $idxTabor = $i
    }
    if ($txt.StartsWith("A csoportve")) {
        $idxCsoport = $i
    }
    if ($txt.StartsWith("Az ifjúsági")) {
        $idxIfjusagi = $i
    }
}

$content8 = $bulletPPr + '<w:r><w:t>A táborvezető funkciók osztálydiagramjai –</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Lovas István, Balla Gábor</w:t></w:r>'
Replace-ParaContent $idxTabor $content8

$content9 = $bulletPPr + '<w:r><w:t>A csoportvezető funkciók osztálydiagramjai –</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Lovas István, Balla Gábor</w:t></w:r>'
Replace-ParaContent $idxCsoport $content9

$content10 = $bulletPPr + '<w:r><w:t>Az ifjúsági vezető funkciók osztálydiagramjai –</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Lovas István, Balla Gábor</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Replace-ParaContent $idxIfjusagi $content10

Write-Host "Done."
